$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add 10 new rows of institution test-data (DEC_0116 .. DEC_0125), following
# exactly the same pattern as the existing rows above (21-32).
# ---------------------------------------------------------------------------
$constUsuario    = "13712759-8"
$constPassword   = "Verity1.0"
$constBusqueda   = "verity"
$constPrefijo    = "ASDF"
$constDescrip    = "Prueba Automation QA"
$constEmail      = "pruebaAutomationQA@acepta.com"
$constNombreRep  = "Nombre de Prueba"

$startRow = 33
$startNum = 116

for ($i = 0; $i -lt 10; $i++) {
    $row = $startRow + $i
    $num = $startNum + $i
    $dec = "DEC_0" + $num
    $empresa = "Nueva Empresa QA " + $num

    $ws.Cells.Item($row, 1).Value  = $dec            # A  TC
    $ws.Cells.Item($row, 2).Value  = $constUsuario    # B  USUARIO
    $ws.Cells.Item($row, 3).Value  = $constPassword   # C  PASSWORD
    $ws.Cells.Item($row, 4).Value  = $constBusqueda   # D  BUSQUEDA
    $ws.Cells.Item($row, 5).Value  = $empresa          # E  NOMBRE_INSTITUCION
    $ws.Cells.Item($row, 6).Value  = $constPrefijo    # F  PREFIJO
    $ws.Cells.Item($row, 7).Value  = $constDescrip    # G  DESCRIPCION
    $ws.Cells.Item($row, 8).Value  = $constEmail      # H  EMAIL
    $ws.Cells.Item($row, 9).Value  = $constNombreRep  # I  NOMBRE_REPRESENTANTE
    $ws.Cells.Item($row, 10).Value = $constEmail      # J  EMAIL_REPRESENTANTE
}

Write-Output "rows written"

# ---------------------------------------------------------------------------
# Turn the new H and J column "EMAIL" cells into real mailto: hyperlinks,
# exactly like the existing rows above. Hyperlinks.Add() forces Excel's
# built-in "Hyperlink" cell style (blue/underline) onto the target cell, so
# immediately re-apply the original cell formatting (copied from the
# unaffected neighbouring row) right after creating each link.
# ---------------------------------------------------------------------------

# Column J (rows 33-42) first, then column H (rows 33-42) -- matching the
# creation order used when the source workbook was produced.
for ($i = 0; $i -lt 10; $i++) {
    $row = $startRow + $i
    $target = $ws.Range("J$row")
    $ws.Hyperlinks.Add($target, "mailto:$constEmail") | Out-Null
    $ws.Range("J32").Copy()
    $target.PasteSpecial(-4122) | Out-Null
}

for ($i = 0; $i -lt 10; $i++) {
    $row = $startRow + $i
    $target = $ws.Range("H$row")
    $ws.Hyperlinks.Add($target, "mailto:$constEmail") | Out-Null
    $ws.Range("H32").Copy()
    $target.PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = $false

Write-Output "hyperlinks added"
